# Update cryptocurrency price/volume data on the "cryptos" worksheet.
# Values are plain text cells (t="inlineStr" in the original file); numeric-looking
# "Price" strings are prefixed with a leading apostrophe so Excel keeps them as text
# instead of silently converting them to floating point numbers (which would drop
# trailing zeros / change their representation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.095.68'
$ws.Range("E2").Value = '  -0.57%  '

$ws.Range("D3").Value = '2.917.04'
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = "'356.87"
$ws.Range("E5").Value = '  +1.18%  '

$ws.Range("D6").Value = "'109.95"
$ws.Range("E6").Value = '  -2.02%  '

$ws.Range("E7").Value = '  +1.82%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").Value = "'0.634"
$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("D10").Value = "'38.94"

$ws.Range("E11").Value = '  +1.27%  '

$ws.Range("D12").Value = "'0.0870"
$ws.Range("E12").Value = '  +0.50%  '

$ws.Range("D13").Value = "'19.53"
$ws.Range("E13").Value = '  -1.74%  '

$ws.Range("E14").Value = '  -0.29%  '

$ws.Range("D15").Value = '3.376.28'
$ws.Range("E15").Value = '  +0.17%  '

$ws.Range("D16").Value = '2.897.48'
$ws.Range("E16").Value = '  -0.55%  '

$ws.Range("D17").Value = "'0.987"
$ws.Range("E17").Value = '  -2.08%  '

$ws.Range("D18").Value = '52.065.56'
$ws.Range("E18").Value = '  -0.67%  '

$ws.Range("D19").Value = "'3.47"
$ws.Range("E19").Value = '  +4.51%  '

$ws.Range("D20").Value = "'7.55"
$ws.Range("E20").Value = '  -1.25%  '

$ws.Range("D21").Value = "'13.91"
$ws.Range("E21").Value = '  -1.84%  '

$ws.Range("E22").Value = '  +0.19%  '

$ws.Range("D23").Value = "'70.60"
$ws.Range("E23").Value = '  -0.50%  '

$ws.Range("D24").Value = "'268.58"
$ws.Range("E24").Value = '  -0.61%  '

$ws.Range("D25").Value = "'2.82"
$ws.Range("E25").Value = '  +1.50%  '

$ws.Range("E26").Value = '  +8.57%  '

$ws.Range("D27").Value = "'7.69"
$ws.Range("E27").Value = '  +16.29%  '

$ws.Range("D28").Value = "'26.95"
$ws.Range("E28").Value = '  +0.69%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("E30").Value = '  +7.46%  '

$ws.Range("E31").Value = '  -1.52%  '

$ws.Range("D32").Value = "'37.40"
$ws.Range("E32").Value = '  -1.25%  '

$ws.Range("D33").Value = "'6.16"
$ws.Range("E33").Value = '  -3.38%  '

$ws.Range("D34").Value = "'2.21"
$ws.Range("E34").Value = '  -1.84%  '

$ws.Range("D35").Value = "'52.17"
$ws.Range("E35").Value = '  -2.47%  '

$ws.Range("D36").Value = "'0.0442"
$ws.Range("E36").Value = '  -1.82%  '

$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("E38").Value = '  -3.72%  '

$ws.Range("D39").Value = "'18.27"
$ws.Range("E39").Value = '  -2.87%  '

$ws.Range("E40").Value = '  -3.69%  '

$ws.Range("D41").Value = "'2.72"
$ws.Range("E41").Value = '  -4.61%  '

$ws.Range("E42").Value = '  +2.52%  '

$ws.Range("D43").Value = "'22.96"
$ws.Range("E43").Value = '  -2.99%  '

$ws.Range("D44").Value = "'119.73"
$ws.Range("E44").Value = '  -0.91%  '

$ws.Range("E45").Value = '  -0.88%  '

$ws.Range("E46").Value = '  -2.29%  '

$ws.Range("E47").Value = '  -4.85%  '

$ws.Range("D48").Value = '2.128.61'
$ws.Range("E48").Value = '  -3.22%  '

$ws.Range("E49").Value = '  -4.54%  '

$ws.Range("D50").Value = "'0.0345"
$ws.Range("E50").Value = '  +0.73%  '

$ws.Range("D51").Value = "'0.924"
$ws.Range("E51").Value = '  -4.78%  '
